$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the stray entries that were accidentally duplicated into rows 60 and 61
# (date + station name), leaving only the blank, pre-formatted C/D/E cells —
# matching the pattern already used by rows 62/63.
$ws.Range("A60:B61").Clear()

# Update the saved selection to reflect the new last-used cell
$ws.Range("C62").Select()
